# Generate Report for Handoff
#
# - Flip the "Handed back: in sync with en-US" status (Overview!E2:F2,
#   zh-cn!C2, de-de!C2) over to "Ready for handoff".
# - Bump the two "Latest …Datetime" timestamps that accompany the new
#   handoff status.
# - Narrow the now-shorter status columns (Overview!E:F, zh-cn!C, de-de!C)
#   to match their new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps that move forward alongside the new handoff
$wsOverview.Range("G2").Value = "2016-08-22 13:02:41"
$wsDeDe.Range("H2").Value = "2016-08-22 13:02:41"
$wsZhCn.Range("H2").Value = "2016-08-22 13:02:36"

# --- Shrink the status columns now that "Ready for handoff" is shorter
# than "Handed back: in sync with en-US" (the nearest width the engine's
# pixel-snapped ColumnWidth grid can represent for the target ~17.22 chars).
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
